$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04941833333333333
$ws.Range("H2").Value = 0.148255
$ws.Range("I2").Value = 0.005167549122999764
$ws.Range("J2").Value = 0.005167549122999764
$ws.Range("M2").Value = 0.7521946666666667
$ws.Range("N2").Value = 2.256584
$ws.Range("O2").Value = 0.07361670343069449
$ws.Range("P2").Value = 0.0736167034306945
$ws.Range("Q2").Value = 0.03717220676888889
$ws.Range("R2").Value = 0.33454986092
$ws.Range("S2").Value = 0.000380417931251419
$ws.Range("T2").Value = 0.0003804179312514191
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04941833333333333
$ws.Range("H3").Value = 0.148255
$ws.Range("I3").Value = 0.005167549122999764
$ws.Range("J3").Value = 0.005167549122999764
$ws.Range("O3").Value = 0.6908862423022597
$ws.Range("P3").Value = 0.6908862423022598
$ws.Range("Q3").Value = 0.3488578686061111
$ws.Range("R3").Value = 3.139720817455
$ws.Range("S3").Value = 0.003570188595501644
$ws.Range("T3").Value = 0.003570188595501645
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04941833333333333
$ws.Range("H4").Value = 0.148255
$ws.Range("I4").Value = 0.005167549122999764
$ws.Range("J4").Value = 0.005167549122999764
$ws.Range("M4").Value = 2.406242333333334
$ws.Range("N4").Value = 7.218727
$ws.Range("O4").Value = 0.2354970542670457
$ws.Range("P4").Value = 0.2354970542670457
$ws.Range("Q4").Value = 0.1189124857094445
$ws.Range("R4").Value = 1.070212371385
$ws.Range("S4").Value = 0.0012169425962467
$ws.Range("T4").Value = 0.0012169425962467
$ws.Range("I5").Value = 0.806706161560336
$ws.Range("J5").Value = 0.806706161560336
$ws.Range("M5").Value = 0.7521946666666667
$ws.Range("N5").Value = 2.256584
$ws.Range("O5").Value = 0.07361670343069449
$ws.Range("P5").Value = 0.0736167034306945
$ws.Range("Q5").Value = 5.802953687617778
$ws.Range("R5").Value = 52.22658318856001
$ws.Range("S5").Value = 0.05938704825130117
$ws.Range("T5").Value = 0.05938704825130118
$ws.Range("I6").Value = 0.806706161560336
$ws.Range("J6").Value = 0.806706161560336
$ws.Range("O6").Value = 0.6908862423022597
$ws.Range("P6").Value = 0.6908862423022598
$ws.Range("S6").Value = 0.5573421886025002
$ws.Range("T6").Value = 0.5573421886025003
$ws.Range("I7").Value = 0.806706161560336
$ws.Range("J7").Value = 0.806706161560336
$ws.Range("M7").Value = 2.406242333333334
$ws.Range("N7").Value = 7.218727
$ws.Range("O7").Value = 0.2354970542670457
$ws.Range("P7").Value = 0.2354970542670457
$ws.Range("S7").Value = 0.1899769247065346
$ws.Range("T7").Value = 0.1899769247065346
$ws.Range("I8").Value = 0.1881262893166642
$ws.Range("J8").Value = 0.1881262893166643
$ws.Range("M8").Value = 0.7521946666666667
$ws.Range("N8").Value = 2.256584
$ws.Range("O8").Value = 0.07361670343069449
$ws.Range("P8").Value = 0.0736167034306945
$ws.Range("Q8").Value = 1.353266153584889
$ws.Range("R8").Value = 12.179395382264
$ws.Range("S8").Value = 0.0138492372481419
$ws.Range("T8").Value = 0.0138492372481419
$ws.Range("I9").Value = 0.1881262893166642
$ws.Range("J9").Value = 0.1881262893166643
$ws.Range("O9").Value = 0.6908862423022597
$ws.Range("P9").Value = 0.6908862423022598
$ws.Range("S9").Value = 0.1299738651042579
$ws.Range("T9").Value = 0.1299738651042579
$ws.Range("I10").Value = 0.1881262893166642
$ws.Range("J10").Value = 0.1881262893166643
$ws.Range("M10").Value = 2.406242333333334
$ws.Range("N10").Value = 7.218727
$ws.Range("O10").Value = 0.2354970542670457
$ws.Range("P10").Value = 0.2354970542670457
$ws.Range("Q10").Value = 4.329047321557445
$ws.Range("S10").Value = 0.04430318696426441
$ws.Range("T10").Value = 0.04430318696426442
